$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text (no Excel auto-number coercion),
# while keeping the cells style/format identical to before (no left-over "@" format).
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "42.621.50"
Set-TextValue "E2" "  +0.31%  "

Set-TextValue "D3" "2.299.35"
Set-TextValue "E3" "  +0.68%  "

Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.07%  "

Set-TextValue "D5" "301.57"
Set-TextValue "E5" "  -1.23%  "

Set-TextValue "D6" "95.79"
Set-TextValue "E6" "  +0.33%  "

Set-TextValue "D7" "0.505"
Set-TextValue "E7" "  +0.46%  "

Set-TextValue "E8" "  +0.04%  "

Set-TextValue "D9" "0.493"
Set-TextValue "E9" "  +0.01%  "

Set-TextValue "D10" "34.46"
Set-TextValue "E10" "  -0.90%  "

Set-TextValue "D11" "19.18"
Set-TextValue "E11" "  +7.11%  "

Set-TextValue "E12" "  +0.13%  "

Set-TextValue "E13" "  +0.78%  "

Set-TextValue "D14" "6.69"
Set-TextValue "E14" "  +0.62%  "

Set-TextValue "D15" "2.658.63"
Set-TextValue "E15" "  +0.43%  "

Set-TextValue "D16" "2.295.29"
Set-TextValue "E16" "  +0.85%  "

Set-TextValue "E17" "  +1.04%  "

Set-TextValue "D18" "42.540.75"
Set-TextValue "E18" "  +0.27%  "

Set-TextValue "D19" "12.41"
Set-TextValue "E19" "  -1.63%  "

Set-TextValue "D20" "0.0₃0888"
Set-TextValue "E20" "  +0.15%  "

Set-TextValue "D21" "6.03"
Set-TextValue "E21" "  +0.85%  "

Set-TextValue "D22" "67.70"
Set-TextValue "E22" "  +1.18%  "

Set-TextValue "D23" "2.32"
Set-TextValue "E23" "  +7.84%  "

Set-TextValue "D24" "235.85"
Set-TextValue "E24" "  +0.39%  "

Set-TextValue "E25" "  -0.01%  "

Set-TextValue "D26" "2.40"
Set-TextValue "E26" "  -0.98%  "

Set-TextValue "D27" "24.26"
Set-TextValue "E27" "  -2.20%  "

Set-TextValue "E28" "  +15.45%  "

Set-TextValue "D29" "165.36"
Set-TextValue "E29" "  -0.29%  "

Set-TextValue "D30" "9.07"
Set-TextValue "E30" "  +1.36%  "

Set-TextValue "D31" "32.58"
Set-TextValue "E31" "  +0.81%  "

Set-TextValue "E32" "  -0.05%  "

Set-TextValue "E33" "  +1.81%  "

Set-TextValue "D34" "17.64"
Set-TextValue "E34" "  +0.99%  "

Set-TextValue "D35" "4.43"
Set-TextValue "E35" "  -3.86%  "

Set-TextValue "D36" "0.0696"
Set-TextValue "E36" "  +1.95%  "

Set-TextValue "D37" "2.33"
Set-TextValue "E37" "  -1.12%  "

Set-TextValue "D38" "0.0995"
Set-TextValue "E38" "  -0.74%  "

Set-TextValue "E39" "  +0.41%  "

Set-TextValue "E40" "  -0.36%  "

Set-TextValue "E41" "  +1.24%  "

Set-TextValue "D42" "20.10"
Set-TextValue "E42" "  +12.71%  "

Set-TextValue "D43" "1.949.77"
Set-TextValue "E43" "  -2.20%  "

Set-TextValue "E44" "  +0.75%  "

Set-TextValue "D45" "10.32"
Set-TextValue "E45" "  +3.70%  "

Set-TextValue "E46" "  +3.09%  "

Set-TextValue "D47" "2.73"
Set-TextValue "E47" "  -0.70%  "

Set-TextValue "D48" "2.527.12"
Set-TextValue "E48" "  +0.67%  "

Set-TextValue "D49" "53.43"
Set-TextValue "E49" "  +0.36%  "

Set-TextValue "B50" "HuobiToken"
Set-TextValue "C50" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D50" "2.81"
Set-TextValue "E50" "  -2.38%  "

Set-TextValue "B51" "BitcoinSV"
Set-TextValue "C51" "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D51" "71.67"
Set-TextValue "E51" "  +1.20%  "
